$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log rows appended by the SONICViewer run on 2018.08.27
$data = @(
    @("2018.08.27", "16:23:13", "RS", 32, 0, 500, 100, 250, 0.1,   0.9400000000000001, "sonic", 6000, 6.6,   61,  37.45, 0.4808171216961022),
    @("2018.08.27", "16:25:54", "RS", 32, 0, 500, 100, 250, 0.1,   0.96,               "sonic", 6000, 9.32,  85,  36.8,  0.4878273340644669),
    @("2018.08.27", "17:21:14", "RS", 32, 0, 500, 100, 250, "N/A", 1,                  "sonic", 6000, 10.84, 114, 35.8,  0.5267293749593411),
    @("2018.08.27", "17:24:37", "RS", 32, 0, 500, 100, 250, 0.1,   0.25,               "sonic", 6000, 7.22,  0,   "N/A", "N/A"),
    @("2018.08.27", "17:31:04", "RS", 32, 0, 500, 100, 250, 0.1,   0.9500000000000001, "sonic", 6000, 8.84,  78,  37.1,  0.4857362089357296)
)

$startRow = 21
$endRow = $startRow + $data.Length - 1

# Column A holds strings that look like dates ("2018.08.27"); force the
# range to Text first so the engine doesn't coerce them into date serials,
# then strip the formatting back off once the literal strings are in place
# so the new rows end up unstyled, matching the rest of the log.
$dateRange = $ws.Range("A$startRow`:A$endRow")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

$dateRange.ClearFormats()
